$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Proyektor Epson "
$ws.Range("B4").Value = "Kursi"
$ws.Range("B5").Value = "Meja"
$ws.Range("B2").Value = "AC Panasonic"

$ws.Range("A6").Select()
